# Insert a new weekly price record at row 143 for
# "Macroferia Regional de Talca" / Arándano (blue), pushing the
# existing rows 143-154 down to 144-155.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 143 (shifts 143:154 -> 144:155)
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new record
$ws.Cells.Item(143, 1).Value = 5
$ws.Cells.Item(143, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(143, 3).Value = 'Maule'
$ws.Cells.Item(143, 4).Value = 45267
$ws.Cells.Item(143, 5).Value = 7
$ws.Cells.Item(143, 6).Value = 'Fruta'
$ws.Cells.Item(143, 7).Value = 100101
$ws.Cells.Item(143, 8).Value = 'Berries'
$ws.Cells.Item(143, 9).Value = 100101001
$ws.Cells.Item(143, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(143, 11).Value = 'Sin especificar'
$ws.Cells.Item(143, 12).Value = 'Primera'
$ws.Cells.Item(143, 13).Value = 190
$ws.Cells.Item(143, 14).Value = 4000
$ws.Cells.Item(143, 15).Value = 4000
$ws.Cells.Item(143, 16).Value = 4000
$ws.Cells.Item(143, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(143, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(143, 19).Value = 2000
$ws.Cells.Item(143, 20).Value = 2
